$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 ("time_taken"), copying the header style/format
# from E1 (bold, bordered, centered) so it matches the other header cells.
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Populate the new time_taken column for each data row.
$ws.Range("F2").Value = "2021-10-05 13:40:56.733122"
$ws.Range("F3").Value = "2021-10-05 13:40:56.733132"
$ws.Range("F4").Value = "2021-10-05 13:40:56.733135"
$ws.Range("F5").Value = "2021-10-05 13:40:56.733137"
$ws.Range("F6").Value = "2021-10-05 13:40:56.733140"
$ws.Range("F7").Value = "2021-10-05 13:40:56.733142"
$ws.Range("F8").Value = "2021-10-05 13:40:56.733144"
$ws.Range("F9").Value = "2021-10-05 13:40:56.733146"
$ws.Range("F10").Value = "2021-10-05 13:40:56.733149"
